$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 5286.615
$ws.Range("J112").Value = 5560.5415
$ws.Range("L112").Value = 16681.6245
$ws.Range("N112").Value = -18897.6245
$ws.Range("H115").Value = 1343.2273
$ws.Range("I115").Value = 515.5
$ws.Range("J115").Value = 2336.5
$ws.Range("K115").Value = 1546.5
$ws.Range("L115").Value = 7009.5
$ws.Range("M115").Value = 20.5
$ws.Range("N115").Value = -10143.5
$ws.Range("H138").Value = 9358.416999999999
$ws.Range("I138").Value = 7015.8335
$ws.Range("J138").Value = 9826.933999999999
$ws.Range("K138").Value = 21047.5005
$ws.Range("L138").Value = 29480.802
$ws.Range("M138").Value = -15907.5005
$ws.Range("N138").Value = -39760.802

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2845.4783
$ws.Range("I45").Value = 1213.8
$ws.Range("J45").Value = 5904.875
$ws.Range("K45").Value = 1213.8
$ws.Range("L45").Value = 5904.875
$ws.Range("M45").Value = -836.8
$ws.Range("N45").Value = -6658.875
$ws.Range("H61").Value = 3824.984
$ws.Range("I61").Value = 3235.673
$ws.Range("K61").Value = 3235.673
$ws.Range("M61").Value = -3023.673
$ws.Range("H74").Value = 178628.06
$ws.Range("I74").Value = 208753.6
$ws.Range("J74").Value = 3900
$ws.Range("K74").Value = 208753.6
$ws.Range("L74").Value = 3900
$ws.Range("M74").Value = -207879.6
$ws.Range("N74").Value = -5648
$ws.Range("H77").Value = 178628.06
$ws.Range("I77").Value = 208753.6
$ws.Range("J77").Value = 3900
$ws.Range("K77").Value = 1043768
$ws.Range("L77").Value = 19500
$ws.Range("M77").Value = -1039400
$ws.Range("N77").Value = -28236
$ws.Range("H132").Value = 3602.5173
$ws.Range("I132").Value = 2867.3684
$ws.Range("K132").Value = 8602.1052
$ws.Range("M132").Value = -6072.1052
$ws.Range("H136").Value = 3824.984
$ws.Range("I136").Value = 3235.673
$ws.Range("K136").Value = 9707.019
$ws.Range("M136").Value = -7157.019

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 99989.5
$ws.Range("J59").Value = 99989.5
$ws.Range("L59").Value = 99989.5
$ws.Range("N59").Value = -101683.5
$ws.Range("H68").Value = 35000
$ws.Range("J68").Value = 35000
$ws.Range("L68").Value = 35000
$ws.Range("N68").Value = -36622
$ws.Range("H71").Value = 35000
$ws.Range("J71").Value = 35000
$ws.Range("L71").Value = 105000
$ws.Range("N71").Value = -113112

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5808.7915
$ws.Range("I31").Value = 3803.4583
$ws.Range("J31").Value = 7814.125
$ws.Range("K31").Value = 3803.4583
$ws.Range("L31").Value = 7814.125
$ws.Range("M31").Value = -3508.4583
$ws.Range("N31").Value = -8404.125
$ws.Range("H34").Value = 5808.7915
$ws.Range("I34").Value = 3803.4583
$ws.Range("J34").Value = 7814.125
$ws.Range("K34").Value = 3803.4583
$ws.Range("L34").Value = 7814.125
$ws.Range("M34").Value = -3601.4583
$ws.Range("N34").Value = -8218.125
$ws.Range("H62").Value = 13532
$ws.Range("I62").Value = 6749.6665
$ws.Range("K62").Value = 6749.6665
$ws.Range("M62").Value = -6125.6665
$ws.Range("H65").Value = 13532
$ws.Range("I65").Value = 6749.6665
$ws.Range("K65").Value = 33748.3325
$ws.Range("M65").Value = -30628.3325
$ws.Range("H75").Value = 35000
$ws.Range("J75").Value = 35000
$ws.Range("L75").Value = 35000
$ws.Range("N75").Value = -36996
$ws.Range("H78").Value = 35000
$ws.Range("J78").Value = 35000
$ws.Range("L78").Value = 105000
$ws.Range("N78").Value = -114984
$ws.Range("H86").Value = 26063.438
$ws.Range("I86").Value = 31871.05
$ws.Range("K86").Value = 31871.05
$ws.Range("M86").Value = -30748.05
$ws.Range("H89").Value = 26063.438
$ws.Range("I89").Value = 31871.05
$ws.Range("K89").Value = 159355.25
$ws.Range("M89").Value = -153739.25
$ws.Range("H141").Value = 227406.64
$ws.Range("J141").Value = 239130.23
$ws.Range("L141").Value = 239130.23
$ws.Range("N141").Value = -249490.23

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1968.1666
$ws.Range("J45").Value = 1968.1666
$ws.Range("L45").Value = 5904.4998
$ws.Range("N45").Value = -6968.4998
$ws.Range("H81").Value = 994
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 994
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 1683.8572
$ws.Range("I113").Value = 1396.75
$ws.Range("J113").Value = 2066.6667
$ws.Range("K113").Value = 4190.25
$ws.Range("L113").Value = 6200.000100000001
$ws.Range("M113").Value = -2020.25
$ws.Range("N113").Value = -10540.0001
$ws.Range("H122").Value = 2520
$ws.Range("J122").Value = 2678.2
$ws.Range("L122").Value = 24103.8
$ws.Range("N122").Value = -29003.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3260.853
$ws.Range("I126").Value = 3185.625
$ws.Range("K126").Value = 9556.875
$ws.Range("M126").Value = -7086.875
$ws.Range("H131").Value = 89499.5
$ws.Range("J131").Value = 89499.5
$ws.Range("L131").Value = 89499.5
$ws.Range("N131").Value = -99579.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3953
$ws.Range("I7").Value = 3967.2856
$ws.Range("K7").Value = 3967.2856
$ws.Range("M7").Value = -3855.2856
$ws.Range("H40").Value = 8759.15
$ws.Range("I40").Value = 8693.842000000001
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 8693.842000000001
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -8557.842000000001
$ws.Range("N40").Value = -10272
$ws.Range("H46").Value = 2808.9
$ws.Range("I46").Value = 1956
$ws.Range("K46").Value = 1956
$ws.Range("M46").Value = -1768
$ws.Range("H126").Value = 3953
$ws.Range("I126").Value = 3967.2856
$ws.Range("K126").Value = 11901.8568
$ws.Range("M126").Value = -9431.856800000001
$ws.Range("H136").Value = 3443
$ws.Range("I136").Value = 2557.1667
$ws.Range("K136").Value = 7671.500100000001
$ws.Range("M136").Value = -5121.500100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 159074.81
$ws.Range("I126").Value = 2246
$ws.Range("K126").Value = 6738
$ws.Range("M126").Value = -4268
$ws.Range("H132").Value = 5004206.5
$ws.Range("I132").Value = 5956674.5
$ws.Range("J132").Value = 3750.375
$ws.Range("K132").Value = 17870023.5
$ws.Range("L132").Value = 11251.125
$ws.Range("M132").Value = -17867493.5
$ws.Range("N132").Value = -16311.125
